$p = $ppt.ActivePresentation
$s = $p.Slides.Item(5)
$shape = $s.Shapes.Item(2)
$table = $shape.Table
$table.ApplyStyle("{C7514A46-C157-4D84-B604-6D1F9F4F42DE}")
